# Add "Saudi Arabia" as a new country row in each of the three comparison
# groups (Global / High-income / International), inserted right after
# "Russia" and right before "USA" in each group's block of rows.
#
# The workbook lists, for each group, the same 12 "along" categories
# (All, Europe, France, Germany, Italy, Poland, Spain, United Kingdom,
# Switzerland, Japan, Russia, USA) stacked one block after another. We
# insert one new row into each block, immediately before the existing
# "USA" row of that block, and populate it with Saudi Arabia's mean /
# CI_low / CI_high values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$globalLabel = "**Global**:<br>Implemented by<br>All other countries"
$highIncomeLabel = "**High-income**:<br>All other HICs and<br>not some MICs (such as China)"
$internationalLabel = "**International**:<br>Some countries (e.g. EU, UK, Brazil)<br>and not others (e.g. U.S., China)"

# Before the inserts, the "USA" row of each group sits at row 13 (Global),
# row 25 (High-income) and row 37 (International). Insert a blank row
# ahead of each one (processing top-to-bottom, accounting for the shift
# introduced by the previous inserts) and fill it in with Saudi Arabia's
# figures.

# --- Global group: insert before row 13 ---
$ws.Rows(13).Insert()
$ws.Cells.Item(13, 1).Value = "Saudi Arabia"
$ws.Cells.Item(13, 2).Value = 83.8567662269365
$ws.Cells.Item(13, 3).Value = 79.7610338864347
$ws.Cells.Item(13, 4).Value = 87.9524985674384
$ws.Cells.Item(13, 5).Value = $globalLabel

# --- High-income group: insert before what is now row 26 ---
$ws.Rows(26).Insert()
$ws.Cells.Item(26, 1).Value = "Saudi Arabia"
$ws.Cells.Item(26, 2).Value = 83.6736163847433
$ws.Cells.Item(26, 3).Value = 79.7506898985947
$ws.Cells.Item(26, 4).Value = 87.596542870892
$ws.Cells.Item(26, 5).Value = $highIncomeLabel

# --- International group: insert before what is now row 39 ---
$ws.Rows(39).Insert()
$ws.Cells.Item(39, 1).Value = "Saudi Arabia"
$ws.Cells.Item(39, 2).Value = 82.5079657388237
$ws.Cells.Item(39, 3).Value = 78.5050343978621
$ws.Cells.Item(39, 4).Value = 86.5108970797853
$ws.Cells.Item(39, 5).Value = $internationalLabel
